$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Merge the runs/proofErr markers around "las mismas" into a single
#    plain run (text content is unchanged, only run/proofErr splitting).
# ---------------------------------------------------------------------
$find1 = $d.Content.Find
$find1.ClearFormatting()
$find1.Replacement.ClearFormatting()
$find1.Execute(
    " eficiente el consumo de recursos y costos de las mismas, ",
    $false, $false, $false, $false, $false, $true, 1, $false,
    " eficiente el consumo de recursos y costos de las mismas, ",
    2
) | Out-Null

# ---------------------------------------------------------------------
# 2) Merge the runs/proofErr markers around "Brindado resultados" into
#    a single plain run (text content unchanged).
# ---------------------------------------------------------------------
$find2 = $d.Content.Find
$find2.ClearFormatting()
$find2.Replacement.ClearFormatting()
$find2.Execute(
    "corporativos. Brindado resultados por encima",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "corporativos. Brindado resultados por encima",
    2
) | Out-Null

# ---------------------------------------------------------------------
# 3) Append new paragraphs at the end of the document: four blank
#    paragraphs followed by the "reply" text, each separated by a
#    blank paragraph.
# ---------------------------------------------------------------------
function Get-EndRange {
    return $d.Range($d.Content.End - 1, $d.Content.End - 1)
}

$newParas = @(
    "",
    "",
    "",
    "",
    "¡Hola Steve!",
    "",
    "Primero que nada, quiero felicitarte por tu dedicación y esfuerzo en alcanzar tus metas y convertirte en un profesional exitoso.",
    "",
    "En cuanto a tu BMY, me parece que es muy claro y conciso. Destacas tus habilidades en modelamiento e implementación de costos, lo cual es muy importante en el mundo empresarial. También mencionas tu formación como ingeniero economista y analista de sistemas, lo cual te da una perspectiva integral para ofrecer soluciones a las empresas.",
    "",
    "En cuanto a cómo contribuirías a quien te contrate, me parece que es muy valioso que puedas ofrecer soluciones en gestión de costos y brindar herramientas e información para la toma de decisiones.",
    "",
    "En cuanto a cómo te conocen las personas, me parece que es muy interesante que te describas como analítico, exigente y autodidacta. Sin embargo, me gustaría sugerirte que incluyas alguna característica personal que te haga destacar como un profesional único y diferente. Por ejemplo, podrías mencionar tu capacidad para adaptarte a diferentes situaciones o tu habilidad para trabajar bajo presión.",
    "",
    "Espero que estos consejos te sean útiles para mejorar tu BMY. ¡Mucho éxito en tus proyectos y buenos estudios!"
)

foreach ($t in $newParas) {
    $r = Get-EndRange
    $r.InsertParagraphAfter()
    if ($t -ne "") {
        $r2 = Get-EndRange
        $r2.InsertAfter($t)
    }
}
